# Adds the new "ODI Bowling Extra" worksheet (bowling extras scraped data)
# and trims the placeholder-empty cells out of the existing
# "ODI Batting Extra" worksheet, per the additional-scraping commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "ODI Bowling Extra" sheet, positioned right after the
#    existing "ODI Batting Extra" sheet (i.e. as the new last sheet).
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$bowlingExtra = $wb.Worksheets.Add($null, $battingExtra)
$bowlingExtra.Name = "ODI Bowling Extra"

# Reuse the header formatting (bold / border / centered) already used by the
# other sheets' header rows, so the new header row matches the workbook's
# existing look.
$battingExtra.Range("A1:C1").Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122)  # xlPasteFormats

$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $bowlingExtra.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL per match.
# $null entries mean the source data had no value for that match/column.
$data = @(
    @("4522", $null, $null),
    @("4533", $null, $null),
    @("4535", "0", "10.00%"),
    @("4577", "0", "20.00%"),
    @("4580", "0", "40.00%"),
    @("4583", "1", "20.00%"),
    @("4586", $null, $null),
    @("4590", $null, $null),
    @("4592", "0", "10.00%"),
    @("4606", $null, $null),
    @("4611", "0", $null),
    @("4616", "0", $null),
    @("4621", "0", "10.00%"),
    @("4623", $null, $null),
    @("4624", "0", "10.00%"),
    @("4636", "0", "30.00%"),
    @("4639", "1", "20.00%"),
    @("4642", $null, $null),
    @("4727", $null, $null),
    @("4731", "0", "20.00%")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $val = $rowValues[$c]
        if ($val -ne $null) {
            $cell = $bowlingExtra.Cells.Item($rowNum, $c + 1)
            # Force text storage (these are codes/percentages kept as text
            # in the source data, not numbers) instead of letting Excel
            # auto-convert numeric-looking strings to numbers.
            $cell.NumberFormat = "@"
            $cell.Value = $val
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Drop the leftover placeholder-empty cells from "ODI Batting Extra" rows
#    that never had real NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH
#    data scraped for them.
# ---------------------------------------------------------------------------
$clearMap = @{
    2  = @("B", "C", "D", "E");
    3  = @("B", "C", "D", "E");
    6  = @("B", "C", "D", "E");
    7  = @("B", "C", "D", "E");
    9  = @("B", "C", "D", "E");
    13 = @("B", "C", "D", "E");
    15 = @("C", "D", "E");
    17 = @("B", "C", "D", "E");
    18 = @("B", "C", "D", "E");
    20 = @("B", "C", "D", "E", "F");
    21 = @("B", "C", "D", "E", "F")
}

foreach ($rowNum in $clearMap.Keys) {
    foreach ($col in $clearMap[$rowNum]) {
        $battingExtra.Range("$col$rowNum").ClearContents()
    }
}
